# fixed pcfdev install directions
#
# Slide 9 ("Setting up a pcfdev environment") body placeholder: the first
# bullet "Install Vagrant" is replaced with a new first bullet describing
# how to download pcfdev, and "Install Vagrant" becomes its own, second
# bullet.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# First bullet ("Install Vagrant") -> rewritten as the "Download pcfdev ..." bullet.
$para1 = $tr.Paragraphs(1, 1)
$para1.Text = "Download "
$run = $para1.InsertAfter("pcfdev")
$run = $run.InsertAfter(" from https")
$run = $run.InsertAfter("://")
$run = $run.InsertAfter("network.pivotal.io")
$run = $run.InsertAfter("/")

# Split off a brand-new second bullet for "Install Vagrant".
$run = $run.InsertAfter([char]13)
$para2 = $tr.Paragraphs(2, 1)
$para2.Text = "Install "
$null = $para2.InsertAfter("Vagrant")
